# Apply the "Absent" column (H) consolidation:
# H = 1 when Total Attendance Count (D) is 0, otherwise H = 0
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 21; $row++) {
    $dVal = $ws.Cells.Item($row, 4).Value2
    if ($dVal -eq 0) {
        $ws.Cells.Item($row, 8).Value = 1
    } else {
        $ws.Cells.Item($row, 8).Value = 0
    }
}
